# The commit swaps the presentation's applied color theme from the
# "Integral" (Red Violet) palette back to the default "Office" palette,
# and updates three tables (slides 14, 15, 16) so they pick up the
# matching default table style instead of the old one.

$p = $ppt.ActivePresentation

# --- 1. Revert the design/theme color scheme from "Red Violet" back to
#        the stock "Office" colors. The theme is shared across the whole
#        deck, so touching it from any slide updates the master theme. ---
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72

# --- 2. Point the three tables at the new (default) table style id. ---
$newStyleId = "{A7E01C34-90C5-4E30-B23A-9711AC2DC7E7}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    $tableShape.Table.ApplyStyle($newStyleId)
}
